$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("B3").Value = 9394
$ws.Range("C3").Value = 17400
$ws.Range("D3").Value = 30100
$ws.Range("E3").Value = 38500
$ws.Range("F3").Value = 39100
$ws.Range("G3").Value = 39500

# Row 4
$ws.Range("B4").Value = 1231.028224
$ws.Range("C4").Value = 2274.361344
$ws.Range("D4").Value = 3947.88864
$ws.Range("E4").Value = 5041.553408
$ws.Range("F4").Value = 5125.439488
$ws.Range("G4").Value = 5181.014016

# Row 5
$ws.Range("B5").Value = 105.5
$ws.Range("C5").Value = 113.63
$ws.Range("D5").Value = 130.56
$ws.Range("E5").Value = 205.11
$ws.Range("F5").Value = 397.32
$ws.Range("G5").Value = 788.98

# Row 6
$ws.Range("B6").Value = 155
$ws.Range("C6").Value = 182
$ws.Range("D6").Value = 231
$ws.Range("E6").Value = 404
$ws.Range("F6").Value = 1336
$ws.Range("G6").Value = 2868

# Row 7
$ws.Range("B7").Value = 159
$ws.Range("C7").Value = 204
$ws.Range("D7").Value = 262
$ws.Range("E7").Value = 502
$ws.Range("F7").Value = 2024
$ws.Range("G7").Value = 4359

# Row 12
$ws.Range("C12").Value = 71900
$ws.Range("D12").Value = 142000
$ws.Range("E12").Value = 271000
$ws.Range("F12").Value = 482000
$ws.Range("G12").Value = 774000

# Row 13
$ws.Range("C13").Value = 294.649856
$ws.Range("D13").Value = 581.95968
$ws.Range("E13").Value = 1111.49056
$ws.Range("F13").Value = 1973.420032
$ws.Range("G13").Value = 3169.845248

# Row 14
$ws.Range("B14").Value = 27.5332
$ws.Range("C14").Value = 27.40673
$ws.Range("D14").Value = 27.72106
$ws.Range("E14").Value = 28.84652
$ws.Range("F14").Value = 32.19637
$ws.Range("G14").Value = 40.59279

# Row 15
$ws.Range("C15").Value = 44.288
$ws.Range("D15").Value = 45.312
$ws.Range("E15").Value = 48.896
$ws.Range("F15").Value = 59.648
$ws.Range("G15").Value = 82.432

# Row 16
$ws.Range("B16").Value = 45.824
$ws.Range("C16").Value = 46.848
$ws.Range("D16").Value = 53.504
$ws.Range("E16").Value = 63.232
$ws.Range("G16").Value = 102.912

# Row 21
$ws.Range("B21").Value = 12000
$ws.Range("C21").Value = 18100
$ws.Range("D21").Value = 22800
$ws.Range("E21").Value = 26900
$ws.Range("F21").Value = 29300
$ws.Range("G21").Value = 24200

# Row 22
$ws.Range("B22").Value = 1579.155456
$ws.Range("C22").Value = 2375.02464
$ws.Range("D22").Value = 2982.150144
$ws.Range("E22").Value = 3520.069632
$ws.Range("F22").Value = 3834.642432
$ws.Range("G22").Value = 3177.18528

# Row 23
$ws.Range("B23").Value = 49.98
$ws.Range("C23").Value = 51.12
$ws.Range("D23").Value = 55.24
$ws.Range("E23").Value = 79.21
$ws.Range("F23").Value = 190.42
$ws.Range("G23").Value = 574.36

# Row 24
$ws.Range("B24").Value = 70
$ws.Range("C24").Value = 69
$ws.Range("D24").Value = 72
$ws.Range("E24").Value = 106
$ws.Range("F24").Value = 310
$ws.Range("G24").Value = 1893

# Row 25
$ws.Range("B25").Value = 155
$ws.Range("C25").Value = 85
$ws.Range("D25").Value = 103
$ws.Range("E25").Value = 123
$ws.Range("F25").Value = 392

# Row 30
$ws.Range("B30").Value = 239000
$ws.Range("C30").Value = 423000
$ws.Range("D30").Value = 630000
$ws.Range("E30").Value = 771000
$ws.Range("F30").Value = 685000
$ws.Range("G30").Value = 680000

# Row 31
$ws.Range("B31").Value = 979.369984
$ws.Range("C31").Value = 1732.247552
$ws.Range("D31").Value = 2581.594112
$ws.Range("E31").Value = 3158.310912
$ws.Range("F31").Value = 2807.037952
$ws.Range("G31").Value = 2785.017856

# Row 32
$ws.Range("B32").Value = 2.22498
$ws.Range("C32").Value = 2.38901
$ws.Range("D32").Value = 2.39453
$ws.Range("E32").Value = 2.61313
$ws.Range("F32").Value = 4.92057
$ws.Range("G32").Value = 16.83237

# Row 33
$ws.Range("B33").Value = 3.568
$ws.Range("C33").Value = 3.952
$ws.Range("D33").Value = 3.952
$ws.Range("E33").Value = 4.256
$ws.Range("F33").Value = 9.536
$ws.Range("G33").Value = 55.552

# Row 34
$ws.Range("B34").Value = 4.832
$ws.Range("C34").Value = 5.28
$ws.Range("D34").Value = 5.344
$ws.Range("E34").Value = 5.92
$ws.Range("F34").Value = 13.248
$ws.Range("G34").Value = 84.48

# Row 39
$ws.Range("B39").Value = 18300
$ws.Range("C39").Value = 24700
$ws.Range("D39").Value = 24400
$ws.Range("E39").Value = 24700
$ws.Range("F39").Value = 25000
$ws.Range("G39").Value = 25100

# Row 40
$ws.Range("B40").Value = 2397.044736
$ws.Range("C40").Value = 3233.808384
$ws.Range("D40").Value = 3196.059648
$ws.Range("E40").Value = 3244.294144
$ws.Range("F40").Value = 3278.897152
$ws.Range("G40").Value = 3291.480064

# Row 41
$ws.Range("B41").Value = 54.35
$ws.Range("C41").Value = 79.76
$ws.Range("D41").Value = 162.17
$ws.Range("E41").Value = 320.93
$ws.Range("F41").Value = 637.24
$ws.Range("G41").Value = 1270.21

# Row 42
$ws.Range("B42").Value = 92
$ws.Range("C42").Value = 143
$ws.Range("D42").Value = 586
$ws.Range("E42").Value = 1516
$ws.Range("F42").Value = 3752
$ws.Range("G42").Value = 7046

# Row 43
$ws.Range("B43").Value = 95
$ws.Range("C43").Value = 145
$ws.Range("D43").Value = 775
$ws.Range("E43").Value = 2073
$ws.Range("F43").Value = 4555

# Row 48
$ws.Range("B48").Value = 462000
$ws.Range("C48").Value = 790000
$ws.Range("D48").Value = 776000
$ws.Range("E48").Value = 797000
$ws.Range("F48").Value = 802000
$ws.Range("G48").Value = 804000

# Row 49
$ws.Range("B49").Value = 1890.582528
$ws.Range("C49").Value = 3233.808384
$ws.Range("D49").Value = 3177.18528
$ws.Range("E49").Value = 3263.168512
$ws.Range("F49").Value = 3283.091456
$ws.Range("G49").Value = 3291.480064

# Row 50
$ws.Range("B50").Value = 1.94216
$ws.Range("C50").Value = 2.2863
$ws.Range("D50").Value = 4.852270000000001
$ws.Range("E50").Value = 9.76303
$ws.Range("F50").Value = 19.67969
$ws.Range("G50").Value = 39.52358

# Row 51
$ws.Range("B51").Value = 0.908
$ws.Range("C51").Value = 0.972
$ws.Range("D51").Value = 1.064
$ws.Range("E51").Value = 1.192
$ws.Range("F51").Value = 1.256
$ws.Range("G51").Value = 1.368

# Row 52
$ws.Range("B52").Value = 74.24
$ws.Range("C52").Value = 93.696
$ws.Range("D52").Value = 138.24
$ws.Range("E52").Value = 257.024
$ws.Range("F52").Value = 284.672
$ws.Range("G52").Value = 284.672

# Row 57
$ws.Range("B57").Value = 9570
$ws.Range("C57").Value = 13500
$ws.Range("D57").Value = 17100
$ws.Range("E57").Value = 16000
$ws.Range("F57").Value = 16800
$ws.Range("G57").Value = 18900

# Row 58
$ws.Range("B58").Value = 1254.096896
$ws.Range("C58").Value = 1765.801984
$ws.Range("D58").Value = 2236.612608
$ws.Range("E58").Value = 2222.98112
$ws.Range("F58").Value = 2203.058176
$ws.Range("G58").Value = 2478.833664

# Row 59
$ws.Range("B59").Value = 53.12981
$ws.Range("C59").Value = 55.11626
$ws.Range("D59").Value = 56.1767
$ws.Range("E59").Value = 76.49
$ws.Range("F59").Value = 420.63
$ws.Range("G59").Value = 836.58

# Row 60
$ws.Range("B60").Value = 56.064
$ws.Range("C60").Value = 58.112
$ws.Range("D60").Value = 59.648
$ws.Range("E60").Value = 88
$ws.Range("F60").Value = 537
$ws.Range("G60").Value = 898

# Row 61
$ws.Range("B61").Value = 58.112
$ws.Range("C61").Value = 62.208
$ws.Range("D61").Value = 63.232
$ws.Range("E61").Value = 96
$ws.Range("G61").Value = 21365

# Row 66
$ws.Range("B66").Value = 264000
$ws.Range("C66").Value = 388000
$ws.Range("D66").Value = 506000
$ws.Range("E66").Value = 542000
$ws.Range("F66").Value = 605000

# Row 67
$ws.Range("B67").Value = 1082.130432
$ws.Range("C67").Value = 1588.59264
$ws.Range("D67").Value = 2073.034752
$ws.Range("E67").Value = 2218.786816
$ws.Range("F67").Value = 2476.736512
$ws.Range("G67").Value = 2305.818624

# Row 68
$ws.Range("B68").Value = 2.01044
$ws.Range("C68").Value = 2.08489
$ws.Range("D68").Value = 2.12109
$ws.Range("E68").Value = 2.34079
$ws.Range("F68").Value = 4.81608
$ws.Range("G68").Value = 30.10735

# Row 69
$ws.Range("B69").Value = 2.896
$ws.Range("C69").Value = 3.216
$ws.Range("D69").Value = 3.344
$ws.Range("E69").Value = 3.792
$ws.Range("F69").Value = 8.512
$ws.Range("G69").Value = 52

# Row 70
$ws.Range("B70").Value = 4.192
$ws.Range("C70").Value = 4.64
$ws.Range("D70").Value = 4.576000000000001
$ws.Range("E70").Value = 5.792
$ws.Range("F70").Value = 11.968
$ws.Range("G70").Value = 178
